$d = $word.ActiveDocument

# Update the date line in the first paragraph.
# (This runtime's Find.Execute searches the whole document
#  regardless of the Range it is called on, so for the table
#  cells below -- several of which share identical text -- we
#  set Range.Text directly instead of relying on Find/Replace.)
$p = $d.Paragraphs.Item(1)
$p.Range.Text = "2023-11-11 Saturday"

# Update the answer table cells (table 1, 1-based row/col).
$t = $d.Tables.Item(1)

# Cell(1,1): "60÷7=8, 4" -> "34÷4=8, 2"
$c = $t.Cell(1, 1)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "34÷4=8, 2"

# Cell(1,2): "40÷5=8, 0" -> "56÷3=18, 2"
$c = $t.Cell(1, 2)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "56÷3=18, 2"

# Cell(1,3): "68÷3=22, 2" -> "43÷2=21, 1"
$c = $t.Cell(1, 3)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "43÷2=21, 1"

# Cell(1,4): "74÷2=37, 0" -> "95÷5=19, 0"
$c = $t.Cell(1, 4)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "95÷5=19, 0"

# Cell(1,5): "25÷7=3, 4" -> "52÷5=10, 2"
$c = $t.Cell(1, 5)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "52÷5=10, 2"

# Cell(5,1): "80÷8=10, 0" -> "48÷5=9, 3"
$c = $t.Cell(5, 1)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "48÷5=9, 3"

# Cell(5,2): "38÷9=4, 2" -> "12÷5=2, 2"
$c = $t.Cell(5, 2)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "12÷5=2, 2"

# Cell(5,3): "77÷3=25, 2" -> "20÷3=6, 2"
$c = $t.Cell(5, 3)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "20÷3=6, 2"

# Cell(5,4): "71÷8=8, 7" -> "14÷4=3, 2"
$c = $t.Cell(5, 4)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "14÷4=3, 2"

# Cell(5,5): "77÷3=25, 2" -> "65÷4=16, 1"
$c = $t.Cell(5, 5)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "65÷4=16, 1"

# Cell(9,1): "96÷8=12, 0" -> "90÷7=12, 6"
$c = $t.Cell(9, 1)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "90÷7=12, 6"

# Cell(9,2): "52÷4=13, 0" -> "57÷7=8, 1"
$c = $t.Cell(9, 2)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "57÷7=8, 1"

# Cell(9,3): "27÷7=3, 6" -> "80÷3=26, 2"
$c = $t.Cell(9, 3)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "80÷3=26, 2"

# Cell(9,4): "77÷3=25, 2" -> "29÷5=5, 4"
$c = $t.Cell(9, 4)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "29÷5=5, 4"

# Cell(9,5): "54÷2=27, 0" -> "45÷7=6, 3"
$c = $t.Cell(9, 5)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "45÷7=6, 3"

# Cell(13,1): "26÷3=8, 2" -> "19÷2=9, 1"
$c = $t.Cell(13, 1)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "19÷2=9, 1"

# Cell(13,2): "85÷7=12, 1" -> "45÷9=5, 0"
$c = $t.Cell(13, 2)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "45÷9=5, 0"

# Cell(13,3): "17÷2=8, 1" -> "59÷7=8, 3"
$c = $t.Cell(13, 3)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "59÷7=8, 3"

# Cell(13,4): "23÷9=2, 5" -> "26÷9=2, 8"
$c = $t.Cell(13, 4)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "26÷9=2, 8"

# Cell(13,5): "25÷6=4, 1" -> "27÷7=3, 6"
$c = $t.Cell(13, 5)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "27÷7=3, 6"

# Cell(17,1): "84÷8=10, 4" -> "10÷4=2, 2"
$c = $t.Cell(17, 1)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "10÷4=2, 2"

# Cell(17,2): "39÷9=4, 3" -> "74÷2=37, 0"
$c = $t.Cell(17, 2)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "74÷2=37, 0"

# Cell(17,3): "19÷2=9, 1" -> "15÷6=2, 3"
$c = $t.Cell(17, 3)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "15÷6=2, 3"

# Cell(17,4): "60÷5=12, 0" -> "26÷6=4, 2"
$c = $t.Cell(17, 4)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "26÷6=4, 2"

# Cell(17,5): "46÷6=7, 4" -> "57÷3=19, 0"
$c = $t.Cell(17, 5)
$r = $c.Range
$r.MoveEnd(1, -1)
$r.Text = "57÷3=19, 0"
